$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 103
$ws.Range("C103").Value = 99
$ws.Range("D103").Value = 44389
$ws.Range("D103").NumberFormat = "d-mmm-yy"
$ws.Range("E103").Value = "Learing"
$ws.Range("F103").Value = "JBPM api program"
$ws.Range("G103").Value = 3

# Row 104
$ws.Range("C104").Value = 100
$ws.Range("D104").Value = 44389
$ws.Range("D104").NumberFormat = "d-mmm-yy"
$ws.Range("E104").Value = "Training"
$ws.Range("F104").Value = "JBPM meeting for Kie server deploy"
$ws.Range("G104").Value = 5

# Row 105
$ws.Range("C105").Value = 101
$ws.Range("D105").Value = 44389
$ws.Range("D105").NumberFormat = "d-mmm-yy"
$ws.Range("E105").Value = "Learing"
$ws.Range("F105").Value = "JBPM  Interview Qusetion Review"
$ws.Range("G105").Value = 3

# Update the active selection to match the authored state (G105 instead of G106)
[void]$ws.Range("G105").Select()
